# Auto-generated edit script: updates cryptos list (price/volume refresh + 3 row reorderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.177.72"
$ws.Range("E2").Value = "  +0.39%  "

# Row 3
$ws.Range("D3").Value = "1.903.26"
$ws.Range("E3").Value = "  +0.83%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "'305.99"
$ws.Range("E5").Value = "  -0.41%  "

# Row 6
$ws.Range("D6").Value = "'0.9998"

# Row 7
$ws.Range("D7").Value = "'0.5244"
$ws.Range("E7").Value = "  +1.96%  "

# Row 8
$ws.Range("E8").Value = "  +1.13%  "

# Row 9
$ws.Range("D9").Value = "'0.07247"
$ws.Range("E9").Value = "  +0.55%  "

# Row 10
$ws.Range("D10").Value = "'21.16"
$ws.Range("E10").Value = "  +0.44%  "

# Row 11
$ws.Range("D11").Value = "'0.9001"
$ws.Range("E11").Value = "  -0.28%  "

# Row 12
$ws.Range("D12").Value = "'0.08418"
$ws.Range("E12").Value = "  +10.17%  "

# Row 13
$ws.Range("D13").Value = "1.886.85"
$ws.Range("E13").Value = "  -0.13%  "

# Row 14
$ws.Range("D14").Value = "'94.92"
$ws.Range("E14").Value = "  +0.56%  "

# Row 15
$ws.Range("D15").Value = "'5.292"
$ws.Range("E15").Value = "  +0.56%  "

# Row 16
$ws.Range("E16").Value = "  +0.15%  "

# Row 17
$ws.Range("D17").Value = "'0.000008624"
$ws.Range("E17").Value = "  +1.35%  "

# Row 18
$ws.Range("D18").Value = "'14.57"
$ws.Range("E18").Value = "  +1.48%  "

# Row 19
$ws.Range("D19").Value = "'0.9998"
$ws.Range("E19").Value = "  +0.07%  "

# Row 20
$ws.Range("D20").Value = "27.224.11"
$ws.Range("E20").Value = "  +0.49%  "

# Row 21
$ws.Range("D21").Value = "'5.066"
$ws.Range("E21").Value = "  +0.17%  "

# Row 22
$ws.Range("D22").Value = "2.154.54"
$ws.Range("E22").Value = "  +0.78%  "

# Row 23
$ws.Range("E23").Value = "  +0.58%  "

# Row 24
$ws.Range("D24").Value = "'6.430"
$ws.Range("E24").Value = "  +0.20%  "

# Row 25
$ws.Range("D25").Value = "'147.09"
$ws.Range("E25").Value = "  +0.66%  "

# Row 26
$ws.Range("D26").Value = "'2.274"
$ws.Range("E26").Value = "  +4.37%  "

# Row 27
$ws.Range("E27").Value = "  -2.26%  "

# Row 28
$ws.Range("E28").Value = "  +0.74%  "

# Row 29
$ws.Range("D29").Value = "'114.93"
$ws.Range("E29").Value = "  +0.40%  "

# Row 30
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'4.816"
$ws.Range("E30").Value = "  -0.22%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.909"
$ws.Range("E31").Value = "  -1.33%  "

# Row 32
$ws.Range("D32").Value = "'0.09258"
$ws.Range("E32").Value = "  +0.41%  "

# Row 33
$ws.Range("D33").Value = "'0.8084"
$ws.Range("E33").Value = "  +5.30%  "

# Row 34
$ws.Range("D34").Value = "'0.05061"

# Row 35
$ws.Range("D35").Value = "'1.239"
$ws.Range("E35").Value = "  +4.00%  "

# Row 36
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "'3.401"
$ws.Range("E36").Value = "  +3.75%  "

# Row 37
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.947"
$ws.Range("E37").Value = "  -0.71%  "

# Row 38
$ws.Range("D38").Value = "'2.626"
$ws.Range("E38").Value = "  +1.67%  "

# Row 39
$ws.Range("D39").Value = "'0.5734"
$ws.Range("E39").Value = "  +1.53%  "

# Row 40
$ws.Range("D40").Value = "'0.01993"
$ws.Range("E40").Value = "  +0.07%  "

# Row 41
$ws.Range("E41").Value = "  -0.22%  "

# Row 42
$ws.Range("D42").Value = "'6.648"
$ws.Range("E42").Value = "  +1.02%  "

# Row 43
$ws.Range("D43").Value = "'8.973"
$ws.Range("E43").Value = "  +0.43%  "

# Row 44
$ws.Range("D44").Value = "'116.66"
$ws.Range("E44").Value = "  -2.03%  "

# Row 45
$ws.Range("D45").Value = "'0.1512"
$ws.Range("E45").Value = "  +0.74%  "

# Row 46
$ws.Range("D46").Value = "'0.4873"
$ws.Range("E46").Value = "  +1.19%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.20"
$ws.Range("E47").Value = "  +0.21%  "

# Row 48
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "'0.9994"
$ws.Range("E48").Value = "  +0.04%  "

# Row 49
$ws.Range("D49").Value = "'1.614"
$ws.Range("E49").Value = "  +1.66%  "

# Row 50
$ws.Range("E50").Value = "  +0.66%  "

# Row 51
$ws.Range("D51").Value = "'63.98"
$ws.Range("E51").Value = "  +0.53%  "
